$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 53.40575035845904
$ws.Range("B3").Value = 0.943946663536994
$ws.Range("B4").Value = 0.05458037482045512
$ws.Range("B5").Value = 0.3756374367841092
